# "ya casi funciona todo" — tidy up the users/machine workbook:
#  - drop the scratch "machine" sheet (denomination table, no longer needed)
#  - post the missing saldo_corriente / saldo_tarjeta_credito figures for every user
#  - a couple of balance corrections on row 5 (yully velandia) and row 6 (Cristian)
#  - leave the "users" sheet focused with F12 selected, like it was left in Excel

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# --- remove the helper "machine" sheet -------------------------------------
if ($wb.Worksheets.Count -gt 1) {
    $excel.DisplayAlerts = $false
    $wb.Worksheets.Item("machine").Delete() | Out-Null
    $excel.DisplayAlerts = $true
}

# --- cell value corrections on the "users" sheet ----------------------------
$ws.Range("G2").Value = 5000000
$ws.Range("I2").Value = 5000000

$ws.Range("G3").Value = 5000000
$ws.Range("I3").Value = 5000000

$ws.Range("G4").Value = 5000000
$ws.Range("I4").Value = 5000000

$ws.Range("E5").Value = 2949800
$ws.Range("G5").Value = 4750000
$ws.Range("I5").Value = 5000000
$ws.Range("K5").Value = 900000

$ws.Range("E6").Value = 4771000
$ws.Range("G6").Value = 5000000
$ws.Range("I6").Value = 5000000

$ws.Range("G7").Value = 5000000
$ws.Range("I7").Value = 5000000

# --- column widths (best-fit like the columns ended up after autosizing) ---
$ws.Columns.Item(2).ColumnWidth  = 11 - (5/6)
$ws.Columns.Item(3).ColumnWidth  = 10.7109375 - (5/6)
$ws.Columns.Item(4).ColumnWidth  = 14.7109375 - (5/6)
$ws.Columns.Item(5).ColumnWidth  = 13.42578125 - (5/6)
$ws.Columns.Item(6).ColumnWidth  = 16.28515625 - (5/6)
$ws.Columns.Item(7).ColumnWidth  = 14.85546875 - (5/6)
$ws.Columns.Item(8).ColumnWidth  = 14.140625 - (5/6)
$ws.Columns.Item(9).ColumnWidth  = 20 - (5/6)
$ws.Columns.Item(10).ColumnWidth = 23.140625 - (5/6)
$ws.Columns.Item(11).ColumnWidth = 18.7109375 - (5/6)
$ws.Columns.Item(12).ColumnWidth = 11.140625 - (5/6)

# --- leave "users" as the active sheet with F12 selected --------------------
$ws.Activate()
$ws.Range("F12").Select() | Out-Null
